$wb = $excel.ActiveWorkbook

# ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 902.2
$ws.Range("I2").Value = 307.75
$ws.Range("J2").Value = 1793.875
$ws.Range("K2").Value = 307.75
$ws.Range("L2").Value = 1793.875
$ws.Range("M2").Value = -194.75
$ws.Range("N2").Value = -2019.875
$ws.Range("H19").Value = 1074.0769
$ws.Range("I19").Value = 986.8125
$ws.Range("J19").Value = 1213.7
$ws.Range("K19").Value = 986.8125
$ws.Range("L19").Value = 1213.7
$ws.Range("M19").Value = -811.8125
$ws.Range("N19").Value = -1563.7
$ws.Range("H51").Value = 11862.5
$ws.Range("I51").Value = 9000
$ws.Range("J51").Value = 12816.667
$ws.Range("K51").Value = 9000
$ws.Range("L51").Value = 12816.667
$ws.Range("M51").Value = -8516
$ws.Range("N51").Value = -13784.667
# ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 7189.5693
$ws.Range("I32").Value = 3111.1755
$ws.Range("J32").Value = 36248.125
$ws.Range("K32").Value = 3111.1755
$ws.Range("L32").Value = 36248.125
$ws.Range("M32").Value = -2824.1755
$ws.Range("N32").Value = -36822.125
$ws.Range("H45").Value = 2196.5
$ws.Range("I45").Value = 999.2
$ws.Range("J45").Value = 4192
$ws.Range("K45").Value = 999.2
$ws.Range("L45").Value = 4192
$ws.Range("M45").Value = -622.2
$ws.Range("N45").Value = -4946
$ws.Range("H61").Value = 4779.8423
$ws.Range("I61").Value = 4779.8423
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4779.8423
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4778.1113
$ws.Range("N61").ClearContents()
$ws.Range("H110").Value = 2809.7693
$ws.Range("I110").Value = 3051.2727
$ws.Range("J110").Value = 1481.5
$ws.Range("K110").Value = 3051.2727
$ws.Range("L110").Value = 1481.5
$ws.Range("M110").Value = -1006.2727
$ws.Range("N110").Value = -5571.5
$ws.Range("H122").Value = 3893.3333
$ws.Range("I122").Value = 3641.8635
$ws.Range("J122").Value = 4999.8
$ws.Range("K122").Value = 10925.5905
$ws.Range("L122").Value = 14999.4
$ws.Range("M122").Value = -8475.5905
$ws.Range("N122").Value = -19899.4
$ws.Range("H136").Value = 4779.8423
$ws.Range("I136").Value = 4779.8423
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 14339.5269
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -12420.3339
$ws.Range("N136").ClearContents()
# BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 39988
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 119964
$ws.Range("N78").ClearContents()
$ws.Range("H82").Value = 90545.11
$ws.Range("I82").Value = 2475
$ws.Range("J82").Value = 115708
$ws.Range("K82").Value = 2475
$ws.Range("L82").Value = 115708
$ws.Range("M82").Value = -2092
$ws.Range("N82").Value = -116474
$ws.Range("H85").Value = 90545.11
$ws.Range("I85").Value = 2475
$ws.Range("J85").Value = 115708
$ws.Range("K85").Value = 2475
$ws.Range("L85").Value = 115708
$ws.Range("M85").Value = -1149
$ws.Range("N85").Value = -118360
$ws.Range("H99").Value = 1469.5
$ws.Range("I99").Value = 1299.5
$ws.Range("J99").Value = 2999.5
$ws.Range("K99").Value = 1299.5
$ws.Range("L99").Value = 2999.5
$ws.Range("M99").Value = 198.5
$ws.Range("N99").Value = -5995.5
$ws.Range("H105").Value = 3924.389
$ws.Range("I105").Value = 3807.8667
$ws.Range("J105").Value = 4507
$ws.Range("K105").Value = 3807.8667
$ws.Range("L105").Value = 4507
$ws.Range("M105").Value = -2060.8667
$ws.Range("N105").Value = -8001
$ws.Range("H134").Value = 2306.9678
$ws.Range("I134").Value = 1672.4828
$ws.Range("J134").Value = 11507
$ws.Range("K134").Value = 5017.4484
$ws.Range("L134").Value = 34521
$ws.Range("M134").Value = -2482.4484
$ws.Range("N134").Value = -39591
# CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 9373
$ws.Range("I31").Value = 3123.6924
$ws.Range("J31").Value = 17497.1
$ws.Range("K31").Value = 3123.6924
$ws.Range("L31").Value = 17497.1
$ws.Range("M31").Value = -2828.6924
$ws.Range("N31").Value = -18087.1
$ws.Range("H34").Value = 9373
$ws.Range("I34").Value = 3123.6924
$ws.Range("J34").Value = 17497.1
$ws.Range("K34").Value = 3123.6924
$ws.Range("L34").Value = 17497.1
$ws.Range("M34").Value = -2921.6924
$ws.Range("N34").Value = -17901.1
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 80000
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 240000
$ws.Range("N78").ClearContents()
$ws.Range("H134").Value = 9093064
$ws.Range("I134").Value = 9616760
$ws.Range("J134").Value = 15665.667
$ws.Range("K134").Value = 28850280
$ws.Range("L134").Value = 46997.001
$ws.Range("M134").Value = -28847745
$ws.Range("N134").Value = -52067.001
# CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H76").Value = 11829.909
$ws.Range("I76").Value = 11654
$ws.Range("J76").Value = 12137.75
$ws.Range("K76").Value = 34962
$ws.Range("L76").Value = 36413.25
$ws.Range("M76").Value = -34579
$ws.Range("N76").Value = -37179.25
$ws.Range("H79").Value = 11829.909
$ws.Range("I79").Value = 11654
$ws.Range("J79").Value = 12137.75
$ws.Range("K79").Value = 34962
$ws.Range("L79").Value = 36413.25
$ws.Range("M79").Value = -33636
$ws.Range("N79").Value = -39065.25
# GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H97").Value = 907.1818
$ws.Range("I97").Value = 562.4
$ws.Range("J97").Value = 1194.5
$ws.Range("K97").Value = 562.4
$ws.Range("L97").Value = 1194.5
$ws.Range("M97").Value = -66.39999999999998
$ws.Range("N97").Value = -2186.5
$ws.Range("H107").Value = 650.1667
$ws.Range("I107").Value = 580.4
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 580.4
$ws.Range("L107").Value = 999
$ws.Range("M107").Value = 1339.6
$ws.Range("N107").Value = -4839
# LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 5143.1763
$ws.Range("I7").Value = 5283.0625
$ws.Range("J7").Value = 2905
$ws.Range("K7").Value = 5283.0625
$ws.Range("L7").Value = 2905
$ws.Range("M7").Value = -5171.0625
$ws.Range("N7").Value = -3129
$ws.Range("H100").Value = 1132.75
$ws.Range("I100").Value = 1001.5
$ws.Range("J100").Value = 1264
$ws.Range("K100").Value = 1001.5
$ws.Range("L100").Value = 1264
$ws.Range("M100").Value = -460.5
$ws.Range("N100").Value = -2346
$ws.Range("H126").Value = 5143.1763
$ws.Range("I126").Value = 5283.0625
$ws.Range("J126").Value = 2905
$ws.Range("K126").Value = 15849.1875
$ws.Range("L126").Value = 8715
$ws.Range("M126").Value = -13379.1875
$ws.Range("N126").Value = -13655
$ws.Range("H136").Value = 5497.7812
$ws.Range("I136").Value = 2778.6
$ws.Range("J136").Value = 7897.0586
$ws.Range("K136").Value = 8335.799999999999
$ws.Range("L136").Value = 23691.1758
$ws.Range("M136").Value = -5785.799999999999
$ws.Range("N136").Value = -28791.1758
# WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value = 3685.5173
$ws.Range("I132").Value = 2509.2273
$ws.Range("J132").Value = 7382.4287
$ws.Range("K132").Value = 7527.6819
$ws.Range("L132").Value = 22147.2861
$ws.Range("M132").Value = -4997.6819
$ws.Range("N132").Value = -27207.2861
